$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.575558
$ws.Range("H2").Value = 10.726674
$ws.Range("I2").Value = 0.025194653521236
$ws.Range("J2").Value = 0.02519465352123599
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 520.9614631988574
$ws.Range("R2").Value = 4688.653168789716
$ws.Range("S2").Value = 0.007220630862734733
$ws.Range("T2").Value = 0.007220630862734732
$ws.Range("G3").Value = 3.575558
$ws.Range("H3").Value = 10.726674
$ws.Range("I3").Value = 0.025194653521236
$ws.Range("J3").Value = 0.02519465352123599
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 603.5531272674214
$ws.Range("R3").Value = 5431.978145406792
$ws.Range("S3").Value = 0.008365367970382278
$ws.Range("T3").Value = 0.008365367970382277
$ws.Range("G4").Value = 3.575558
$ws.Range("H4").Value = 10.726674
$ws.Range("I4").Value = 0.025194653521236
$ws.Range("J4").Value = 0.02519465352123599
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 458.1224973276373
$ws.Range("R4").Value = 4123.102475948735
$ws.Range("S4").Value = 0.006349670132614723
$ws.Range("T4").Value = 0.006349670132614722
$ws.Range("G5").Value = 3.575558
$ws.Range("H5").Value = 10.726674
$ws.Range("I5").Value = 0.025194653521236
$ws.Range("J5").Value = 0.02519465352123599
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 235.1325521070819
$ws.Range("R5").Value = 2116.192968963738
$ws.Range("S5").Value = 0.003258984555504264
$ws.Range("T5").Value = 0.003258984555504263
$ws.Range("I6").Value = 0.7460690747908298
$ws.Range("J6").Value = 0.7460690747908298
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 15426.81412637188
$ws.Range("R6").Value = 138841.3271373469
$ws.Range("S6").Value = 0.2138187525629578
$ws.Range("T6").Value = 0.2138187525629578
$ws.Range("I7").Value = 0.7460690747908298
$ws.Range("J7").Value = 0.7460690747908298
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.247716934733293
$ws.Range("T7").Value = 0.247716934733293
$ws.Range("I8").Value = 0.7460690747908298
$ws.Range("J8").Value = 0.7460690747908298
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 13566.0142114678
$ws.Range("R8").Value = 122094.1279032102
$ws.Range("S8").Value = 0.1880276907588301
$ws.Range("T8").Value = 0.1880276907588301
$ws.Range("I9").Value = 0.7460690747908298
$ws.Range("J9").Value = 0.7460690747908298
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 6962.791746902787
$ws.Range("R9").Value = 62665.12572212509
$ws.Range("S9").Value = 0.09650569673574892
$ws.Range("T9").Value = 0.09650569673574892
$ws.Range("G10").Value = 32.36130266666667
$ws.Range("H10").Value = 97.08390800000001
$ws.Range("I10").Value = 0.2280292497513723
$ws.Range("J10").Value = 0.2280292497513723
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 4715.06589691672
$ws.Range("R10").Value = 42435.59307225049
$ws.Range("S10").Value = 0.06535176349907712
$ws.Range("T10").Value = 0.06535176349907712
$ws.Range("G11").Value = 32.36130266666667
$ws.Range("H11").Value = 97.08390800000001
$ws.Range("I11").Value = 0.2280292497513723
$ws.Range("J11").Value = 0.2280292497513723
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 5462.578268039341
$ws.Range("R11").Value = 49163.20441235408
$ws.Range("S11").Value = 0.07571243559958472
$ws.Range("T11").Value = 0.07571243559958472
$ws.Range("G12").Value = 32.36130266666667
$ws.Range("H12").Value = 97.08390800000001
$ws.Range("I12").Value = 0.2280292497513723
$ws.Range("J12").Value = 0.2280292497513723
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 4146.329270684147
$ws.Range("R12").Value = 37316.96343615732
$ws.Range("S12").Value = 0.05746895924916853
$ws.Range("T12").Value = 0.05746895924916853
$ws.Range("G13").Value = 32.36130266666667
$ws.Range("H13").Value = 97.08390800000001
$ws.Range("I13").Value = 0.2280292497513723
$ws.Range("J13").Value = 0.2280292497513723
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 2128.114181205577
$ws.Range("R13").Value = 19153.0276308502
$ws.Range("S13").Value = 0.02949609140354194
$ws.Range("T13").Value = 0.02949609140354194
$ws.Range("G14").Value = 0.1003386666666667
$ws.Range("H14").Value = 0.301016
$ws.Range("I14").Value = 0.000707021936561918
$ws.Range("J14").Value = 0.0007070219365619179
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 14.61941845219378
$ws.Range("R14").Value = 131.574766069744
$ws.Range("S14").Value = 0.0002026280858145739
$ws.Range("T14").Value = 0.0002026280858145739
$ws.Range("G15").Value = 0.1003386666666667
$ws.Range("H15").Value = 0.301016
$ws.Range("I15").Value = 0.000707021936561918
$ws.Range("J15").Value = 0.0007070219365619179
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 16.93713709930311
$ws.Range("R15").Value = 152.434233893728
$ws.Range("S15").Value = 0.0002347521333241406
$ws.Range("T15").Value = 0.0002347521333241405
$ws.Range("G16").Value = 0.1003386666666667
$ws.Range("H16").Value = 0.301016
$ws.Range("I16").Value = 0.000707021936561918
$ws.Range("J16").Value = 0.0007070219365619179
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 12.85600752438045
$ws.Range("R16").Value = 115.704067719424
$ws.Range("S16").Value = 0.0001781868549970992
$ws.Range("T16").Value = 0.0001781868549970992
$ws.Range("G17").Value = 0.1003386666666667
$ws.Range("H17").Value = 0.301016
$ws.Range("I17").Value = 0.000707021936561918
$ws.Range("J17").Value = 0.0007070219365619179
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 6.598378985421333
$ws.Range("R17").Value = 59.38541086879199
$ws.Range("S17").Value = 0.00009145486242610446
$ws.Range("T17").Value = 0.00009145486242610444
